# Rename "CO2 electrolyzer purchase cost ($/m^2)" column to
# "hydr. electrolyzer purchase cost ($/m^2)" and move it so it sits right
# after "hydr. electrolyzer current density (A/m^2)" (column AI), pushing
# "hydr. separator energy (unit TBD)" one column to the right (column AJ).
# In other words: swap the whole contents (header label, data values, and
# cell comments) currently in columns AI and AJ, and rename what lands in
# AI to the new "hydr. electrolyzer purchase cost ($/m^2)" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture current values before any writes (rows 1-5, cols AI & AJ) ---
$rows = 1..5
$aiVals = @{}
$ajVals = @{}
foreach ($r in $rows) {
    $aiVals[$r] = $ws.Cells.Item($r, 35).Value()
    $ajVals[$r] = $ws.Cells.Item($r, 36).Value()
}

# --- 2. Capture existing comments on AI1 and AJ2 (the only two comments in
#        these two columns) before moving anything ---
$commentAI1 = $null
if ($ws.Range("AI1").Comment -ne $null) {
    $commentAI1 = $ws.Range("AI1").Comment.Text()
}
$commentAJ2 = $null
if ($ws.Range("AJ2").Comment -ne $null) {
    $commentAJ2 = $ws.Range("AJ2").Comment.Text()
}

# Remove the old comments so they don't linger on the wrong cell once the
# values are swapped.
if ($ws.Range("AI1").Comment -ne $null) {
    $ws.Range("AI1").Comment.Delete() | Out-Null
}
if ($ws.Range("AJ2").Comment -ne $null) {
    $ws.Range("AJ2").Comment.Delete() | Out-Null
}

# --- 3. Write the swapped values back: AJ gets AI's old content, AI gets
#        AJ's old content (row 1 header text is handled specially below so
#        it picks up the renamed label) ---
foreach ($r in $rows) {
    if ($r -eq 1) {
        # Header row: AI1 becomes the renamed label, AJ1 becomes the old
        # AI1 label ("hydr. separator energy (unit TBD)").
        $ws.Cells.Item(1, 35).Value = "hydr. electrolyzer purchase cost (`$/m^2)"
        $ws.Cells.Item(1, 36).Value = $aiVals[1]
    } else {
        $ws.Cells.Item($r, 35).Value = $ajVals[$r]
        $ws.Cells.Item($r, 36).Value = $aiVals[$r]
    }
}

# --- 4. Re-create the comments on their new cells ---
if ($commentAI1 -ne $null) {
    $ws.Range("AJ1").AddComment($commentAI1) | Out-Null
}
if ($commentAJ2 -ne $null) {
    $ws.Range("AI2").AddComment($commentAJ2) | Out-Null
}

# --- 5. Swap the column widths of AI and AJ to match the swapped content ---
$widthAI = $ws.Columns.Item(35).ColumnWidth()
$widthAJ = $ws.Columns.Item(36).ColumnWidth()
$ws.Columns.Item(35).ColumnWidth = $widthAJ
$ws.Columns.Item(36).ColumnWidth = $widthAI

# --- 6. Update the active selection to AI1, matching the new focal point
#        of the edit ---
$ws.Range("AI1").Select() | Out-Null
